$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 0, 6.201049113329182)
    3 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 0, 5.553084769722144)
    4 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 0, 6.201049113329182)
    5 = @(0.3048080303191223, 0.3127903958511391, 26.21740644021617, 8.660232485948974, 1, 35.49523735233541)
    6 = @(0.127881588408715, 0.04240448674262143, 3.900430680208489, 0.496779210170732, 0, 4.567495965530558)
    7 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 0, 5.553084769722144)
    8 = @(0.6753301551942219, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 0, 3.645393585217082)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
